$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to reflect repulled data
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -6
$ws.Range("F5").Value = -1
$ws.Range("F7").Value = -8
$ws.Range("F8").Value = -1
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = -7
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -1
